# Applies the cryptos.xlsx data refresh described in the diff.
# Plain text/percent cells are set directly. Numeric-looking "Price"
# strings (e.g. "143.83") must stay text (matching the source's
# inlineStr cells) instead of being auto-coerced to floating point
# numbers, so for those we briefly force a text NumberFormat, set
# the value, then restore the original cell style to avoid leaving
# stray style artifacts behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.355.09"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "2.451.50"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  -0.09%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.38"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +1.22%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.83"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "2.446.75"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("E11").Value = "  +2.43%  "
$ws.Range("E12").Value = "  -0.73%  "
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.346"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  -3.01%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.31"
$ws.Range("D14").Style = $origStyle
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000175"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "62.135.16"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "2.444.92"
$ws.Range("E18").Value = "  +0.71%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.89"
$ws.Range("D19").Style = $origStyle
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.12"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -2.31%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.29"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("E22").Value = "  -1.50%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.95"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -7.88%  "
$ws.Range("E24").Value = "  +0.05%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.60"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +0.38%  "
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.30"
$ws.Range("D26").Style = $origStyle
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "591.76"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -5.55%  "
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("D29").Value = "0.0₃0958"
$ws.Range("E29").Value = "  -4.45%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("E31").Value = "  -4.29%  "
$ws.Range("E32").Value = "  -2.01%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.89"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  -1.10%  "
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.93"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  -4.12%  "
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("E37").Value = "  -4.17%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.379"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +0.20%  "
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "152.57"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  +4.39%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.34"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -0.24%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.42"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  -2.23%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "43.04"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.72"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -2.64%  "
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -4.79%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.65"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -3.29%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.97"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -3.44%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.607"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  +1.15%  "
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0523"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -1.54%  "
$ws.Range("D50").Value = "0.0₆0246"
$ws.Range("E50").Value = "  +6.83%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.76"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -5.09%  "
